$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the "Resolving-Mac" target-cluster row) entirely.
$ws.Rows("5:5").Delete()

# Update recomputed TPM-derived values in row 2 (FAPs -> ECs)
$ws.Range("M2").Value2 = 30.58864766666666
$ws.Range("N2").Value2 = 91.76594299999999
$ws.Range("O2").Value2 = 0.3925391465174898
$ws.Range("P2").Value2 = 0.3925391465174898
$ws.Range("Q2").Value2 = 4.693889161745333
$ws.Range("R2").Value2 = 42.245002455708
$ws.Range("S2").Value2 = 0.3925391465174898
$ws.Range("T2").Value2 = 0.3925391465174898

# Update recomputed TPM-derived values in row 3 (FAPs -> FAPs)
$ws.Range("O3").Value2 = 0.291183949679193
$ws.Range("P3").Value2 = 0.291183949679193
$ws.Range("S3").Value2 = 0.291183949679193
$ws.Range("T3").Value2 = 0.291183949679193

# Update recomputed TPM-derived values in row 4 (FAPs -> MuSCs)
$ws.Range("M4").Value2 = 24.64590566666666
$ws.Range("N4").Value2 = 73.93771699999999
$ws.Range("O4").Value2 = 0.3162769038033173
$ws.Range("P4").Value2 = 0.3162769038033172
$ws.Range("Q4").Value2 = 3.781963516361333
$ws.Range("R4").Value2 = 34.037671647252
$ws.Range("S4").Value2 = 0.3162769038033173
$ws.Range("T4").Value2 = 0.3162769038033172
